$wb = $excel.ActiveWorkbook

$wsNd = $wb.Worksheets.Item("Neodymium")
$wsNd.Range("C2").Value = 0.000002735394444618311
$wsNd.Range("D2").Value = 0.007285353541089507
$wsNd.Range("E2").Value = 0.008339924611837421
$wsNd.Range("B3").Value = 0.000000000003676322720692059
$wsNd.Range("C3").Value = 0.000132372580562949
$wsNd.Range("D3").Value = 0.006352624942659592
$wsNd.Range("E3").Value = 0.007415669097852164
$wsNd.Range("B4").Value = 0.00000000000005738846834589288
$wsNd.Range("C4").Value = 0.0001197210825341109
$wsNd.Range("D4").Value = 0.005191851368327277
$wsNd.Range("E4").Value = 0.006548198374626368
$wsNd.Range("C5").Value = 0.000000002664856196090389
$wsNd.Range("D5").Value = 0.00028701895446782
$wsNd.Range("E5").Value = 0.0005468359805677819

$wsCu = $wb.Worksheets.Item("Copper")
$wsCu.Range("B2").Value = 0.000006274887434311204
$wsCu.Range("C2").Value = 0.004837162278782569
$wsCu.Range("D2").Value = 0.6980914194676436
$wsCu.Range("E2").Value = 0.6373369448677954
$wsCu.Range("B3").Value = 0.00004266705571568352
$wsCu.Range("C3").Value = 0.01745069221883031
$wsCu.Range("D3").Value = 0.4957277151421089
$wsCu.Range("E3").Value = 0.4893987397993288
$wsCu.Range("B4").Value = 0.0001265372631661763
$wsCu.Range("C4").Value = 0.00467118311447194
$wsCu.Range("D4").Value = 0.4207052905870581
$wsCu.Range("E4").Value = 0.4929022936717287
$wsCu.Range("B5").Value = 0.00003975125710766059
$wsCu.Range("C5").Value = 0.01024166722731971
$wsCu.Range("D5").Value = 0.604205072935894
$wsCu.Range("E5").Value = 0.497817373632984

$wsSi = $wb.Worksheets.Item("Raw silicon")
$wsSi.Range("B2").Value = 0.00000107623096286272
$wsSi.Range("C2").Value = 0.00007486642629472632
$wsSi.Range("D2").Value = 0.02016567546390563
$wsSi.Range("E2").Value = 0.01868247443550294
$wsSi.Range("B3").Value = 0.000001148540841621533
$wsSi.Range("C3").Value = 0.0002501457324738229
$wsSi.Range("D3").Value = 0.01062027883849227
$wsSi.Range("E3").Value = 0.01038192857648591
$wsSi.Range("B4").Value = 0.000007359536893654886
$wsSi.Range("C4").Value = 0.00007021125419683742
$wsSi.Range("D4").Value = 0.01097468629610411
$wsSi.Range("E4").Value = 0.01303346366569932
$wsSi.Range("B5").Value = 0.000003951808733664171
$wsSi.Range("C5").Value = 0.00008916107934033047
$wsSi.Range("D5").Value = 0.01879760438469067
$wsSi.Range("E5").Value = 0.01549935507400286

$wb.Save()
